$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the grid-on-photo "completed" mark (row 8 -> "grid on photo") in column C,
# matching the style/value already used by the other "Completed" cells in column C.
$ws.Range("C8").Value = "y"

# Reflect the new active selection on the sheet (as left by the edit).
$ws.Range("C8").Select()
